$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared-string-backed cell values for row 2 ---
# Order matters: shared strings are appended in the order values are first set,
# and the target sharedStrings order is: EFT facade text (idx 9), SQL text (idx 10), Harisha (idx 11)

# C2: short descriptive text
$ws.Range("C2").Value = "EFT Changes. Façade, Service Layer."

# G2: the big SQL / deployment script text
$sql = @"
CREATE TABLE MOLSAEFTUSERCONFIGURATION(CONFIGURATIONID BIGINT not null, USERNAME1 CHARACTER(64), USERNAME2 CHARACTER(64), 
USERTITLE1 CHARACTER(100), USERTITLE2 CHARACTER(100), EFFECTIVEDATETIME DATE, RECORDSTATUS CHARACTER(10), VERSIONNO INT not null, LASTWRITTEN DATE);
ALTER TABLE MOLSAEFTUSERCONFIGURATION ADD CONSTRAINT MOLSAEFTUSERCONFIGURATION PRIMARY KEY(CONFIGURATIONID);
INSERT INTO FunctionIdentifier (fidName, projectPackage, codePackage, description, fidEnabled) values ('MOLSAEFTUserConfiguration.configureUsersForEFTLetterConfiguration', 'curam', 'molsa.eft.eftletter.facade', 'curam.molsa.eft.eftletter.facade.MOLSAEFTUserConfiguration.configureUsersForEFTLetterConfiguration', 'Y');
INSERT INTO FunctionIdentifier (fidName, projectPackage, codePackage, description, fidEnabled) values ('MOLSAEFTUserConfiguration.listAllUsersByPosition', 'curam', 'molsa.eft.eftletter.facade', 'curam.molsa.eft.eftletter.facade.MOLSAEFTUserConfiguration.listAllUsersByPosition', 'Y');
INSERT INTO FunctionIdentifier (fidName, projectPackage, codePackage, description, fidEnabled) values ('MOLSAEFTUserConfiguration.listUsersConfiguredForEFTLetter', 'curam', 'molsa.eft.eftletter.facade', 'curam.molsa.eft.eftletter.facade.MOLSAEFTUserConfiguration.listUsersConfiguredForEFTLetter', 'Y');
INSERT INTO SECURITYIDENTIFIER (DESCRIPTION, LASTWRITTEN, SIDNAME, SIDTYPE, VERSIONNO) VALUES (null, '2015-03-03 03:42:25', 'MOLSAEFTUserConfiguration.configureUsersForEFTLetterConfiguration', 'FUNCTION', 0);
INSERT INTO SECURITYIDENTIFIER (DESCRIPTION, LASTWRITTEN, SIDNAME, SIDTYPE, VERSIONNO) VALUES (null, '2015-03-03 03:42:25', 'MOLSAEFTUserConfiguration.listAllUsersByPosition', 'FUNCTION', 0);
INSERT INTO SECURITYIDENTIFIER (DESCRIPTION, LASTWRITTEN, SIDNAME, SIDTYPE, VERSIONNO) VALUES (null, '2015-03-03 03:42:25', 'MOLSAEFTUserConfiguration.listUsersConfiguredForEFTLetter', 'FUNCTION', 0);
INSERT INTO SECURITYFIDSID(SIDNAME, FIDNAME, LASTWRITTEN) VALUES ('MOLSAEFTUserConfiguration.configureUsersForEFTLetterConfiguration','MOLSAEFTUserConfiguration.configureUsersForEFTLetterConfiguration');
INSERT INTO SECURITYFIDSID(SIDNAME, FIDNAME, LASTWRITTEN) VALUES ('MOLSAEFTUserConfiguration.listAllUsersByPosition','MOLSAEFTUserConfiguration.listAllUsersByPosition');
INSERT INTO SECURITYFIDSID(SIDNAME, FIDNAME, LASTWRITTEN) VALUES('MOLSAEFTUserConfiguration.listUsersConfiguredForEFTLetter','MOLSAEFTUserConfiguration.listUsersConfiguredForEFTLetter');
INSERT INTO SECURITYGROUPSID (GROUPNAME, LASTWRITTEN, SIDNAME) VALUES ('SUPERGROUP', null, 'MOLSAEFTUserConfiguration.configureUsersForEFTLetterConfiguration');
INSERT INTO SECURITYGROUPSID (GROUPNAME, LASTWRITTEN, SIDNAME) VALUES ('SUPERGROUP', null, 'MOLSAEFTUserConfiguration.listAllUsersByPosition');
INSERT INTO SECURITYGROUPSID (GROUPNAME, LASTWRITTEN, SIDNAME) VALUES ('SUPERGROUP', null, 'MOLSAEFTUserConfiguration.listUsersConfiguredForEFTLetter');
--Manger job Update
update job set NAME='مدير ادارة الضمان الاجتماعي' , comments='مدير ادارة الضمان الاجتماعي' where jobid=45014;
--Need to update Assistance Manager Job
INSERT INTO JOB (COMMENTS, JOBID, NAME, RECORDSTATUS, VERSIONNO) VALUES ('Assistance manager', 45021, 'Assistance manager', 'RST1', 1);
update POSITION  set jobid=45021 where positionid=45290;
--New Job for General Secrartey وكيل الوزارة المساعد للشؤون الاجتماعية
INSERT INTO JOB (COMMENTS, JOBID, NAME, RECORDSTATUS, VERSIONNO) VALUES ('General  Secretary', 45020, 'General Secretary', 'RST1', 1);
INSERT INTO POSITION (COMMENTS, FROMDATE, JOBID, LASTWRITTEN, LEADPOSITIONIND, NAME, POSITIONID, RECORDSTATUS, TODATE, VERSIONNO) VALUES ('General Secretary', '2003-01-01', 45020, '2006-01-01 00:00:00', '0', 'General Secretary', 45349, 'RST1', null, 1);
INSERT INTO ORGUNITPOSITIONLINK (ORGUNITPOSITIONLINKID, ORGANISATIONUNITID, POSITIONID, ORGANISATIONSTRUCTUREID, RECORDSTATUS, VERSIONNO) VALUES (45356, 45002, 45349,45000, 'RST1', 1);

"@
$ws.Range("G2").Value = $sql

# E2: developer name
$ws.Range("E2").Value = "Harisha"

# --- Column G width (grew substantially to host the large SQL text) ---
$ws.Columns.Item(7).ColumnWidth = 187

# --- Row 2 height grows to the Excel maximum to accommodate the wrapped SQL text ---
$ws.Rows.Item(2).RowHeight = 409.5

# --- Selection / view state ---
$ws.Range("G2").Select()
